$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.120.38'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '1.564.60'
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.84'
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.491'
$ws.Range("E6").Value = '  +0.99%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.14'
$ws.Range("E8").Value = '  +3.85%  '
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0860'
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("D12").Value = '1.789.40'
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("D13").Value = '1.566.23'
$ws.Range("E13").Value = '  +1.60%  '
$ws.Range("E14").Value = '  +2.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("E15").Value = '  +2.18%  '
$ws.Range("D16").Value = '27.128.18'
$ws.Range("E16").Value = '  +1.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.00'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.77'
$ws.Range("E18").Value = '  +2.57%  '
$ws.Range("D19").Value = '0.0₃0697'
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.36'
$ws.Range("E20").Value = '  +1.62%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.07'
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.95'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.94'
$ws.Range("E25").Value = '  +1.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.63'
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.00'
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  +1.51%  '
$ws.Range("E30").Value = '  +2.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.10'
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("D33").Value = '1.446.20'
$ws.Range("E33").Value = '  +5.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  +4.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.965'
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0166'
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.522'
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("E40").Value = '  +1.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.73'
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.28'
$ws.Range("E43").Value = '  +3.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.986'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.48'
$ws.Range("E45").Value = '  +2.36%  '
$ws.Range("E46").Value = '  +2.49%  '
$ws.Range("D47").Value = '1.702.55'
$ws.Range("E47").Value = '  +1.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.92'
$ws.Range("E48").Value = '  +3.12%  '
$ws.Range("E49").Value = '  +3.05%  '
$ws.Range("D50").Value = '0.0₆0101'
$ws.Range("E50").Value = '  +3.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0968'
$ws.Range("E51").Value = '  +2.74%  '
